# Applies the crypto price/volume refresh described in the commit:
# "Updated cryptos list on Wed May 15 16:41:03 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text that can look like a plain decimal (e.g. "581.42").
# Force the whole data range to Text first so Excel does not silently coerce
# those assignments into numbers (which would also eat significant trailing zeros).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "64.810.75"
$ws.Range("E2").Value = "  +5.35%  "

$ws.Range("D3").Value = "2.972.19"
$ws.Range("E3").Value = "  +2.72%  "

$ws.Range("E4").Value = "  +0.18%  "

$ws.Range("D5").Value = "581.42"
$ws.Range("E5").Value = "  +1.93%  "

$ws.Range("D6").Value = "152.03"
$ws.Range("E6").Value = "  +6.29%  "

$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("D8").Value = "0.514"
$ws.Range("E8").Value = "  +1.33%  "

$ws.Range("D9").Value = "2.970.09"
$ws.Range("E9").Value = "  +2.68%  "

$ws.Range("D10").Value = "6.95"
$ws.Range("E10").Value = "  +3.33%  "

$ws.Range("D11").Value = "0.150"
$ws.Range("E11").Value = "  +2.33%  "

$ws.Range("D12").Value = "0.446"
$ws.Range("E12").Value = "  +2.78%  "

$ws.Range("E13").Value = "  +1.58%  "

$ws.Range("D14").Value = "33.97"
$ws.Range("E14").Value = "  +6.35%  "

$ws.Range("E15").Value = "  +0.75%  "

$ws.Range("D16").Value = "64.807.43"
$ws.Range("E16").Value = "  +5.34%  "

$ws.Range("D17").Value = "3.467.48"
$ws.Range("E17").Value = "  +2.79%  "

$ws.Range("D18").Value = "6.86"
$ws.Range("E18").Value = "  +3.48%  "

$ws.Range("D19").Value = "2.977.33"
$ws.Range("E19").Value = "  +3.56%  "

$ws.Range("D20").Value = "446.07"
$ws.Range("E20").Value = "  +3.08%  "

$ws.Range("D21").Value = "13.60"
$ws.Range("E21").Value = "  +3.23%  "

$ws.Range("D22").Value = "0.676"
$ws.Range("E22").Value = "  +3.15%  "

$ws.Range("E23").Value = "  +4.79%  "

$ws.Range("D24").Value = "80.82"
$ws.Range("E24").Value = "  +1.70%  "

$ws.Range("D25").Value = "12.22"
$ws.Range("E25").Value = "  +3.56%  "

$ws.Range("D26").Value = "10.56"
$ws.Range("E26").Value = "  +4.99%  "

$ws.Range("D27").Value = "2.17"
$ws.Range("E27").Value = "  +7.32%  "

$ws.Range("E28").Value = "  -0.10%  "

$ws.Range("D29").Value = "7.76"
$ws.Range("E29").Value = "  +10.72%  "

$ws.Range("D30").Value = "2.33"
$ws.Range("E30").Value = "  +13.03%  "

$ws.Range("E31").Value = "  +2.52%  "

$ws.Range("E32").Value = "  -2.24%  "

$ws.Range("E33").Value = "  +2.94%  "

$ws.Range("D34").Value = "26.49"
$ws.Range("E34").Value = "  +3.43%  "

$ws.Range("E35").Value = "  +0.07%  "

$ws.Range("E36").Value = "  +2.23%  "

$ws.Range("E37").Value = "  +4.02%  "

$ws.Range("D38").Value = "2.07"
$ws.Range("E38").Value = "  +6.79%  "

$ws.Range("D39").Value = "48.86"
$ws.Range("E39").Value = "  +0.13%  "

$ws.Range("D40").Value = "2.87"
$ws.Range("E40").Value = "  +1.54%  "

$ws.Range("D41").Value = "43.39"
$ws.Range("E41").Value = "  +9.42%  "

$ws.Range("E42").Value = "  +4.11%  "

$ws.Range("D43").Value = "0.293"
$ws.Range("E43").Value = "  +9.63%  "

$ws.Range("E44").Value = "  +1.37%  "

$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "2.760.60"
$ws.Range("E45").Value = "  +2.78%  "

$ws.Range("B46").Value = "Bittensor"
$ws.Range("C46").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D46").Value = "376.85"
$ws.Range("E46").Value = "  +10.83%  "

$ws.Range("D47").Value = "0.0346"
$ws.Range("E47").Value = "  +3.63%  "

$ws.Range("D48").Value = "133.85"
$ws.Range("E48").Value = "  +0.33%  "

$ws.Range("E49").Value = "  -0.03%  "

$ws.Range("E50").Value = "  +1.84%  "

$ws.Range("D51").Value = "22.85"
$ws.Range("E51").Value = "  +5.92%  "

# Restore the default (unstyled) cell style on column D now that the text is set,
# so the saved style table matches the original (no stray explicit number format).
$ws.Range("D2:D51").Style = "Normal"
